# Apply cryptos list update (values scraped from coinranking.com)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep a Text format so that numeric-looking
# values (e.g. "22.20", "0.999") are preserved exactly as strings,
# matching the original inline-string cell contents.
$editRanges = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D27","E27","D28","E28","E29","D30","E30","D31","E31","D32","E32","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","E40","D41","E41","E42","D44","E44","D45","E45","D46","E46","E47","B48","C48","D48","E48","B49","C49","D49","E49","D50","E50","E51")
foreach ($addr in $editRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '77.308.45'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '3.141.80'
$ws.Range('E3').Value = '  +5.86%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '202.54'
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('D6').Value = '630.39'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.227'
$ws.Range('E8').Value = '  +13.99%  '
$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  +4.55%  '
$ws.Range('D10').Value = '3.140.10'
$ws.Range('E10').Value = '  +5.80%  '
$ws.Range('D11').Value = '0.535'
$ws.Range('E11').Value = '  +24.14%  '
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').Value = '5.43'
$ws.Range('E13').Value = '  +9.32%  '
$ws.Range('D14').Value = '3.717.80'
$ws.Range('E14').Value = '  +5.81%  '
$ws.Range('D15').Value = '0.0000223'
$ws.Range('E15').Value = '  +19.41%  '
$ws.Range('D16').Value = '30.51'
$ws.Range('E16').Value = '  +6.01%  '
$ws.Range('D17').Value = '77.243.86'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '3.132.28'
$ws.Range('E18').Value = '  +5.62%  '
$ws.Range('D20').Value = '9.19'
$ws.Range('E20').Value = '  +5.63%  '
$ws.Range('D21').Value = '429.77'
$ws.Range('E21').Value = '  +15.98%  '
$ws.Range('D22').Value = '2.84'
$ws.Range('E22').Value = '  +26.87%  '
$ws.Range('D23').Value = '4.76'
$ws.Range('E23').Value = '  +11.22%  '
$ws.Range('D24').Value = '6.75'
$ws.Range('E24').Value = '  +5.19%  '
$ws.Range('D25').Value = '3.302.14'
$ws.Range('E25').Value = '  +5.75%  '
$ws.Range('E26').Value = '  +3.93%  '
$ws.Range('D27').Value = '4.64'
$ws.Range('E27').Value = '  +8.79%  '
$ws.Range('D28').Value = '10.64'
$ws.Range('E28').Value = '  +11.28%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '0.0000114'
$ws.Range('E30').Value = '  +8.79%  '
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = '8.72'
$ws.Range('E32').Value = '  +5.82%  '
$ws.Range('E33').Value = '  +7.52%  '
$ws.Range('D34').Value = '519.47'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').Value = '1.96'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '0.136'
$ws.Range('E36').Value = '  +22.86%  '
$ws.Range('D37').Value = '22.20'
$ws.Range('E37').Value = '  +10.06%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').Value = '163.72'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').Value = '194.88'
$ws.Range('E41').Value = '  +7.11%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D44').Value = '0.107'
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('D45').Value = '5.37'
$ws.Range('E45').Value = '  +9.96%  '
$ws.Range('D46').Value = '0.796'
$ws.Range('E46').Value = '  +13.36%  '
$ws.Range('E47').Value = '  +7.26%  '
$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').Value = '1.28'
$ws.Range('E48').Value = '  +6.26%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '42.52'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = '2.51'
$ws.Range('E50').Value = '  +10.18%  '
$ws.Range('E51').Value = '  +6.14%  '
